# Update the Max Weight constraint for GLD (row 8) in the
# basic_asset_classes_constrained workbook from 5% to 1%.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 0.01

# Match the author's resulting active-cell selection (C8).
$ws.Range("C8").Select() | Out-Null
